$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'" + '62.796.74'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  -0.97%  '
$c = $ws.Range('D3')
$c.Value = "'" + '2.444.52'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('E4').Value = '  +0.01%  '
$c = $ws.Range('D5')
$c.Value = "'" + '576.21'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.14%  '
$c = $ws.Range('D6')
$c.Value = "'" + '144.30'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.23%  '
$ws.Range('E7').Value = '  +0.13%  '
$c = $ws.Range('D8')
$c.Value = "'" + '0.531'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -2.09%  '
$c = $ws.Range('D9')
$c.Value = "'" + '2.440.79'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -1.67%  '
$ws.Range('E10').Value = '  -4.06%  '
$ws.Range('E11').Value = '  +0.41%  '
$c = $ws.Range('D12')
$c.Value = "'" + '5.21'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -1.16%  '
$c = $ws.Range('D13')
$c.Value = "'" + '0.351'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -2.03%  '
$c = $ws.Range('D14')
$c.Value = "'" + '26.57'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -2.20%  '
$ws.Range('E15').Value = '  -3.36%  '
$ws.Range('E16').Value = '  -1.67%  '
$c = $ws.Range('D17')
$c.Value = "'" + '62.497.68'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.86%  '
$c = $ws.Range('D18')
$c.Value = "'" + '2.434.61'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -1.24%  '
$ws.Range('E19').Value = '  -3.61%  '
$ws.Range('E20').Value = '  -2.67%  '
$c = $ws.Range('D21')
$c.Value = "'" + '330.23'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('E22').Value = '  -0.93%  '
$ws.Range('E23').Value = '  +3.42%  '
$c = $ws.Range('D24')
$c.Value = "'" + '1.00'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +0.19%  '
$c = $ws.Range('D25')
$c.Value = "'" + '65.74'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.57%  '
$c = $ws.Range('D26')
$c.Value = "'" + '635.95'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +0.38%  '
$c = $ws.Range('D27')
$c.Value = "'" + '9.10'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +6.77%  '
$c = $ws.Range('D28')
$c.Value = "'" + '0.0₃0972'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -6.33%  '
$c = $ws.Range('D30')
$c.Value = "'" + '1.00'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('E31').Value = '  -3.75%  '
$ws.Range('E32').Value = '  -2.69%  '
$ws.Range('E33').Value = '  -1.64%  '
$ws.Range('E34').Value = '  -3.86%  '
$c = $ws.Range('D35')
$c.Value = "'" + '5.01'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -3.12%  '
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('E37').Value = '  -2.33%  '
$ws.Range('E38').Value = '  -2.26%  '
$c = $ws.Range('D39')
$c.Value = "'" + '18.54'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.49%  '
$c = $ws.Range('D40')
$c.Value = "'" + '5.27'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -3.33%  '
$c = $ws.Range('D41')
$c.Value = "'" + '146.79'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.45%  '
$c = $ws.Range('D42')
$c.Value = "'" + '1.74'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -3.34%  '
$c = $ws.Range('D43')
$c.Value = "'" + '42.50'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +1.41%  '
$ws.Range('E44').Value = '  -0.02%  '
$c = $ws.Range('D45')
$c.Value = "'" + '2.52'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -4.42%  '
$c = $ws.Range('D46')
$c.Value = "'" + '145.38'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -2.18%  '
$c = $ws.Range('D47')
$c.Value = "'" + '3.71'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.74%  '
$c = $ws.Range('D48')
$c.Value = "'" + '0.0526'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -3.11%  '
$ws.Range('E49').Value = '  -1.45%  '
$c = $ws.Range('D50')
$c.Value = "'" + '19.81'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -5.49%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range('D51')
$c.Value = "'" + '0.0₆0235'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +6.69%  '
